$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1 title
#    paragraph: <w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>
#    <w:r><w:t>: Read our review ...</w:t></w:r>
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs.First
[void]$titlePara.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)

$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:r/>' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' +
  '<w:r><w:t>: Read our review of Break Bones, a 3-reel, 17-fixed-payline video slot game by Hacksaw Gaming. Play for free and learn about its special features.</w:t></w:r>' +
  '</w:p>'
[void]$metaPara.Range.InsertXML($metaXml)

# ---------------------------------------------------------------------------
# 2) Near the end of the document, remove the duplicate bold
#    "Play Break Bones for Free - Review of Hacksaw Gaming's Slot Game"
#    paragraph entirely (text + its own paragraph mark).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$boldPara = $d.Paragraphs($count - 1)
$lastPara = $d.Paragraphs($count)
$killRange = $d.Range($boldPara.Range.Start, $lastPara.Range.Start)
[void]$killRange.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the text of the final (italic) paragraph with the new prompt
#    text, keeping its run formatting (italic) and leading empty run intact.
# ---------------------------------------------------------------------------
$count2 = $d.Paragraphs.Count
$finalPara = $d.Paragraphs($count2)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End)
$finalRange.Text = "Please create a cartoon-style feature image for Break Bones that showcases a happy Maya warrior with glasses. The image should be eye-catching and draw in potential players, highlighting the game's theme and exciting features. Use bold, bright colors to make the image pop, and consider incorporating elements of the game, such as the Wilds and Scatter symbols, into the design. Make sure the Maya warrior is front and center, looking happy and excited to play the game. Overall, the image should convey a sense of fun and adventure, inviting players to join in on the action and give Break Bones a spin."

Write-Output "ok"
